$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ "B" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 2, The Typology)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251399058"; "I" = "https://doi.org/10.1177/27683605251399058"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 2, The Typology)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 3, Using the Classification Model)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251399056"; "I" = "https://doi.org/10.1177/27683605251399056"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 3, Using the Classification Model)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "What is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part I, Conceptual Foundations)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251398692"; "I" = "https://doi.org/10.1177/27683605251398692"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "What is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part I, Conceptual Foundations)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 2, The Typology)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251399058"; "I" = "https://doi.org/10.1177/27683605251399058"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 2, The Typology)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 3, Using the Classification Model)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251399056"; "I" = "https://doi.org/10.1177/27683605251399056"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "What Is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part 3, Using the Classification Model)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "What is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part I, Conceptual Foundations)"; "C" = 2026; "D" = "SAGE Publications"; "E" = "Journal of Integrative and Complementary Medicine"; "F" = "Ijaz, Nadine"; "H" = "10.1177/27683605251398692"; "I" = "https://doi.org/10.1177/27683605251398692"; "J" = "Journal"; "K" = "Co-integration"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "What is Traditional, Complementary, and Integrative Medicine: An Operational Typology (Part I, Conceptual Foundations)"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Nanotechnology"; "F" = "Liu, Wenjun; Zhang, Yachao; Wang, Zhizhe; Su, Kai; Zhao, Shenglei; Xu, Shengrui; Zhang, Jinfeng; Yao, Yixin; Wang, Baiqi; Dong, Yaolong; Hao, Yue; Zhang, Jincheng"; "H" = "10.1088/1361-6528/ae46a5"; "I" = "https://doi.org/10.1088/1361-6528/ae46a5"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Nanotechnology"; "F" = "Liu, Wenjun; Zhang, Yachao; Wang, Zhizhe; Su, Kai; Zhao, Shenglei; Xu, Shengrui; Zhang, Jinfeng; Yao, Yixin; Wang, Baiqi; Dong, Yaolong; Hao, Yue; Zhang, Jincheng"; "H" = "10.1088/1361-6528/ae46a5"; "I" = "https://doi.org/10.1088/1361-6528/ae46a5"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Nanotechnology"; "F" = "Liu, Wenjun; Zhang, Yachao; Wang, Zhizhe; Su, Kai; Zhao, Shenglei; Xu, Shengrui; Zhang, Jinfeng; Yao, Yixin; Wang, Baiqi; Dong, Yaolong; Hao, Yue; Zhang, Jincheng"; "H" = "10.1088/1361-6528/ae46a5"; "I" = "https://doi.org/10.1088/1361-6528/ae46a5"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Journal of Physics D: Applied Physics"; "F" = "Rahman, S. M. Razibur; Rahman, Ehsanur"; "H" = "10.1088/1361-6463/ae46aa"; "I" = "https://doi.org/10.1088/1361-6463/ae46aa"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "C" = 2026; "D" = "Elsevier BV"; "E" = "Solid-State Electronics"; "F" = "Xie, Zijing; Ma, Xiao; Li, Xinghuan; Tang, Jun; Wang, Hong"; "H" = "10.1016/j.sse.2026.109351"; "I" = "https://doi.org/10.1016/j.sse.2026.109351"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Journal of Physics D: Applied Physics"; "F" = "Rahman, S. M. Razibur; Rahman, Ehsanur"; "H" = "10.1088/1361-6463/ae46aa"; "I" = "https://doi.org/10.1088/1361-6463/ae46aa"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "C" = 2026; "D" = "Elsevier BV"; "E" = "Solid-State Electronics"; "F" = "Xie, Zijing; Ma, Xiao; Li, Xinghuan; Tang, Jun; Wang, Hong"; "H" = "10.1016/j.sse.2026.109351"; "I" = "https://doi.org/10.1016/j.sse.2026.109351"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Journal of Physics D: Applied Physics"; "F" = "Rahman, S. M. Razibur; Rahman, Ehsanur"; "H" = "10.1088/1361-6463/ae46aa"; "I" = "https://doi.org/10.1088/1361-6463/ae46aa"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Revealing trap dynamics in p-GaN gate HEMTs: a stretched exponential model for positive and negative bias-temperature instability"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Nanotechnology"; "F" = "Liu, Wenjun; Zhang, Yachao; Wang, Zhizhe; Su, Kai; Zhao, Shenglei; Xu, Shengrui; Zhang, Jinfeng; Yao, Yixin; Wang, Baiqi; Dong, Yaolong; Hao, Yue; Zhang, Jincheng"; "H" = "10.1088/1361-6528/ae46a5"; "I" = "https://doi.org/10.1088/1361-6528/ae46a5"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "C" = 2026; "D" = "IOP Publishing"; "E" = "Nanotechnology"; "F" = "Liu, Wenjun; Zhang, Yachao; Wang, Zhizhe; Su, Kai; Zhao, Shenglei; Xu, Shengrui; Zhang, Jinfeng; Yao, Yixin; Wang, Baiqi; Dong, Yaolong; Hao, Yue; Zhang, Jincheng"; "H" = "10.1088/1361-6528/ae46a5"; "I" = "https://doi.org/10.1088/1361-6528/ae46a5"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Gate Stack"; "Q" = "High-performance GaN HEMTs with over 2 MV/cm breakdown field and 73% PAE via an AlN super back barrier/ultra-thin GaN channel heterostructure"; "R" = "High"; "S" = "2026-02-17" },
    @{ "B" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "C" = 2026; "D" = "Elsevier BV"; "E" = "Solid-State Electronics"; "F" = "Xie, Zijing; Ma, Xiao; Li, Xinghuan; Tang, Jun; Wang, Hong"; "H" = "10.1016/j.sse.2026.109351"; "I" = "https://doi.org/10.1016/j.sse.2026.109351"; "J" = "Journal"; "K" = "n-FET"; "L" = "Experiment"; "M" = "Contacts"; "Q" = "Investigation of cap layer effects on low-contact-resistance vanadium-based Au-free low-temperature ohmic contacts for AlGaN/GaN HEMT"; "R" = "High"; "S" = "2026-02-17" },
)

$startRow = 200
$colMap = @{ "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "Q"=17; "R"=18; "S"=19 }

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $startRow + $i
    foreach ($col in $colMap.Keys) {
        $colNum = $colMap[$col]
        $cell = $ws.Cells.Item($rowNum, $colNum)
        if ($col -eq "S") {
            # Force the date-like text (e.g. "2026-02-17") to be stored as a
            # plain string instead of being auto-converted to a date value,
            # then restore the default "Normal" style so no explicit cell
            # style/number-format index is left behind.
            $cell.Value = "'" + $r[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $r[$col]
        }
    }
}

$wb.Save()
